$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 3-29 with the new data (candidates progressed / new entries)
$ws.Cells.Item(3, 1).Value = 594
$ws.Cells.Item(3, 2).Value = "Antithesis"
$ws.Cells.Item(3, 3).Value = "Enterprise Account Executive"
$ws.Cells.Item(3, 4).Value = "Douglas Hopkins"
$ws.Cells.Item(3, 5).Value = "1st Interview"

$ws.Cells.Item(4, 1).Value = 594
$ws.Cells.Item(4, 2).Value = "Antithesis"
$ws.Cells.Item(4, 3).Value = "Enterprise Account Executive"
$ws.Cells.Item(4, 4).Value = "Edward Aballa"
$ws.Cells.Item(4, 5).Value = "CV Sent"

$ws.Cells.Item(5, 1).Value = 667
$ws.Cells.Item(5, 2).Value = "Antithesis"
$ws.Cells.Item(5, 3).Value = "BDR"
$ws.Cells.Item(5, 4).Value = "Brannon Olive"
$ws.Cells.Item(5, 5).Value = "CV Sent"

$ws.Cells.Item(6, 1).Value = 667
$ws.Cells.Item(6, 2).Value = "Antithesis"
$ws.Cells.Item(6, 3).Value = "BDR"
$ws.Cells.Item(6, 4).Value = "Bromley German"
$ws.Cells.Item(6, 5).Value = "2nd Interview"

$ws.Cells.Item(7, 1).Value = 667
$ws.Cells.Item(7, 2).Value = "Antithesis"
$ws.Cells.Item(7, 3).Value = "BDR"
$ws.Cells.Item(7, 4).Value = "Eray Yaman"
$ws.Cells.Item(7, 5).Value = "1st Interview"

$ws.Cells.Item(8, 1).Value = 667
$ws.Cells.Item(8, 2).Value = "Antithesis"
$ws.Cells.Item(8, 3).Value = "BDR"
$ws.Cells.Item(8, 4).Value = "J. Donahoe"
$ws.Cells.Item(8, 5).Value = "4th Interview"

$ws.Cells.Item(9, 1).Value = 667
$ws.Cells.Item(9, 2).Value = "Antithesis"
$ws.Cells.Item(9, 3).Value = "BDR"
$ws.Cells.Item(9, 4).Value = "Ryan Lewis"
$ws.Cells.Item(9, 5).Value = "4th Interview"

$ws.Cells.Item(10, 1).Value = 680
$ws.Cells.Item(10, 2).Value = "Oscilar"
$ws.Cells.Item(10, 3).Value = "Sales Engineer"
$ws.Cells.Item(10, 4).Value = "JOHN FROST"
$ws.Cells.Item(10, 5).Value = "4th Interview"

$ws.Cells.Item(11, 1).Value = 731
$ws.Cells.Item(11, 2).Value = "Oscilar"
$ws.Cells.Item(11, 3).Value = "Enterprise AE x5"
$ws.Cells.Item(11, 4).Value = "Daniel Grasso"
$ws.Cells.Item(11, 5).Value = "3rd Interview"

$ws.Cells.Item(12, 1).Value = 731
$ws.Cells.Item(12, 2).Value = "Oscilar"
$ws.Cells.Item(12, 3).Value = "Enterprise AE x5"
$ws.Cells.Item(12, 4).Value = "Rob Owen"
$ws.Cells.Item(12, 5).Value = "3rd Interview"

$ws.Cells.Item(13, 1).Value = 731
$ws.Cells.Item(13, 2).Value = "Oscilar"
$ws.Cells.Item(13, 3).Value = "Enterprise AE x5"
$ws.Cells.Item(13, 4).Value = "WILLIAM WOLLISON"
$ws.Cells.Item(13, 5).Value = "2nd Interview"

$ws.Cells.Item(14, 1).Value = 744
$ws.Cells.Item(14, 2).Value = "Synthflow AI"
$ws.Cells.Item(14, 3).Value = "Sales Engineers"
$ws.Cells.Item(14, 4).Value = "Steffen Kaiser"
$ws.Cells.Item(14, 5).Value = "2nd Interview"

$ws.Cells.Item(15, 1).Value = 744
$ws.Cells.Item(15, 2).Value = "Synthflow AI"
$ws.Cells.Item(15, 3).Value = "Sales Engineers"
$ws.Cells.Item(15, 4).Value = "Omer Maroof"
$ws.Cells.Item(15, 5).Value = "2nd Interview"

$ws.Cells.Item(16, 1).Value = 744
$ws.Cells.Item(16, 2).Value = "Synthflow AI"
$ws.Cells.Item(16, 3).Value = "Sales Engineers"
$ws.Cells.Item(16, 4).Value = "John Jänckel"
$ws.Cells.Item(16, 5).Value = "1st Interview"

$ws.Cells.Item(17, 1).Value = 750
$ws.Cells.Item(17, 2).Value = "Novee.io"
$ws.Cells.Item(17, 3).Value = "Enterprise Account Executive (US)"
$ws.Cells.Item(17, 4).Value = "Sean Ribisi"
$ws.Cells.Item(17, 5).Value = "CV Sent"

$ws.Cells.Item(18, 1).Value = 760
$ws.Cells.Item(18, 2).Value = "Impala"
$ws.Cells.Item(18, 3).Value = "Head of Sales (NA)"
$ws.Cells.Item(18, 4).Value = "Patrick Racy"
$ws.Cells.Item(18, 5).Value = "4th Interview"

$ws.Cells.Item(19, 1).Value = 765
$ws.Cells.Item(19, 2).Value = "groundcover"
$ws.Cells.Item(19, 3).Value = "SDR"
$ws.Cells.Item(19, 4).Value = "Rebecca Masters"
$ws.Cells.Item(19, 5).Value = "1st Interview"

$ws.Cells.Item(20, 1).Value = 805
$ws.Cells.Item(20, 2).Value = "Dash0"
$ws.Cells.Item(20, 3).Value = "Dash0 NYC SDR"
$ws.Cells.Item(20, 4).Value = "Rebecca Masters"
$ws.Cells.Item(20, 5).Value = "CV Sent"

$ws.Cells.Item(21, 1).Value = 807
$ws.Cells.Item(21, 2).Value = "Oscilar"
$ws.Cells.Item(21, 3).Value = "SE Leader"
$ws.Cells.Item(21, 4).Value = "Ray Mi"
$ws.Cells.Item(21, 5).Value = "4th Interview"

$ws.Cells.Item(22, 1).Value = 820
$ws.Cells.Item(22, 2).Value = "Silverfort"
$ws.Cells.Item(22, 3).Value = "Nordics RSM"
$ws.Cells.Item(22, 4).Value = "Marc Solis"
$ws.Cells.Item(22, 5).Value = "3rd Interview"

$ws.Cells.Item(23, 1).Value = 820
$ws.Cells.Item(23, 2).Value = "Silverfort"
$ws.Cells.Item(23, 3).Value = "Nordics RSM"
$ws.Cells.Item(23, 4).Value = "Marko Rämö"
$ws.Cells.Item(23, 5).Value = "CV Sent"

$ws.Cells.Item(24, 1).Value = 820
$ws.Cells.Item(24, 2).Value = "Silverfort"
$ws.Cells.Item(24, 3).Value = "Nordics RSM"
$ws.Cells.Item(24, 4).Value = "Peter Inselseth"
$ws.Cells.Item(24, 5).Value = "1st Interview"

$ws.Cells.Item(25, 1).Value = 821
$ws.Cells.Item(25, 2).Value = "Silverfort"
$ws.Cells.Item(25, 3).Value = "DACH - RSM"
$ws.Cells.Item(25, 4).Value = "Jochen Rummel"
$ws.Cells.Item(25, 5).Value = "CV Sent"

$ws.Cells.Item(26, 1).Value = 821
$ws.Cells.Item(26, 2).Value = "Silverfort"
$ws.Cells.Item(26, 3).Value = "DACH - RSM"
$ws.Cells.Item(26, 4).Value = "Narinder Shetge"
$ws.Cells.Item(26, 5).Value = "CV Sent"

$ws.Cells.Item(27, 1).Value = 835
$ws.Cells.Item(27, 2).Value = "Mintlify"
$ws.Cells.Item(27, 3).Value = "TAM"
$ws.Cells.Item(27, 4).Value = "Abbas Engineer"
$ws.Cells.Item(27, 5).Value = "CV Sent"

$ws.Cells.Item(28, 1).Value = 835
$ws.Cells.Item(28, 2).Value = "Mintlify"
$ws.Cells.Item(28, 3).Value = "TAM"
$ws.Cells.Item(28, 4).Value = "Patrick Mackle"
$ws.Cells.Item(28, 5).Value = "1st Interview"

$ws.Cells.Item(29, 1).Value = 835
$ws.Cells.Item(29, 2).Value = "Mintlify"
$ws.Cells.Item(29, 3).Value = "TAM"
$ws.Cells.Item(29, 4).Value = "Vishaal Bhardwaj"
$ws.Cells.Item(29, 5).Value = "1st Interview"

# Remove the now-obsolete trailing rows (30 and 31)
$ws.Rows.Item(31).EntireRow.Delete() | Out-Null
$ws.Rows.Item(30).EntireRow.Delete() | Out-Null
